$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry goes into row 14, right below the last existing history row (13).
# Copy the formatting (number formats / styles) from row 13 so the new row
# matches the date (A) and time (B) formats used by all the other entries.
$ws.Range("A13:D13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Date column (A) - 31/03/2025, stored as Excel serial date 45747
$ws.Range("A14").Value = 45747
# Time column (B) - 10:17 (fraction of a day)
$ws.Range("B14").Value = 0.4284722222222222
# File / version column (C)
$ws.Range("C14").Value = "Futconnect 3103 1017"
# Observation column (D)
$ws.Range("D14").Value = "Implantado filtro de ano no dash e corrigido fórmula da participação."

# After typing into D14, Excel's selection lands on the next empty row below it
$ws.Range("D15").Select()
